$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they stay text, matching the source data
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "51.478.26"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.784.66"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "353.24"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").Value = "108.44"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("D7").Value = "0.551"
$ws.Range("E7").Value = "  -1.59%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  +5.29%  "
$ws.Range("D10").Value = "39.76"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "0.0831"
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("D13").Value = "19.91"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("D15").Value = "3.223.07"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "2.786.63"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "0.939"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "51.397.14"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").Value = "7.67"
$ws.Range("E19").Value = "  +2.87%  "
$ws.Range("D20").Value = "3.17"
$ws.Range("E20").Value = "  +2.93%  "
$ws.Range("D21").Value = "13.37"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "266.79"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").Value = "2.75"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "25.97"
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").Value = "0.164"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "10.30"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "36.69"
$ws.Range("E30").Value = "  +5.29%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "6.17"
$ws.Range("E31").Value = "  +7.10%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "51.86"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "2.08"
$ws.Range("E33").Value = "  -3.28%  "
$ws.Range("D34").Value = "5.65"
$ws.Range("E34").Value = "  +8.64%  "
$ws.Range("D35").Value = "0.0439"
$ws.Range("E35").Value = "  -5.72%  "
$ws.Range("D36").Value = "0.0852"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "18.76"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "3.12"
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").Value = "2.47"
$ws.Range("E42").Value = "  -4.22%  "
$ws.Range("D43").Value = "119.24"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "2.17"
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "21.58"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").Value = "2.117.97"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("E48").Value = "  +5.34%  "
$ws.Range("D49").Value = "1.36"
$ws.Range("E49").Value = "  +8.49%  "
$ws.Range("D50").Value = "0.905"
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("D51").Value = "5.36"
$ws.Range("E51").Value = "  -7.44%  "

# Restore default (unstyled) cell style now that the text values are locked in
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
